# PIWG_A.pptx edit: retitle the slide and strip the leftover "Title #2" demo
# shapes that were pasted under the Technical Recommendation title block.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1. Title: "PIWG Action Item 18-002" -> "PIWG" + " Action Item XX-XXX"
#     (two runs, matching a retype where "PIWG" trips the spell checker)
$title = $s.Shapes.Item("Title 1")
$tr = $title.TextFrame.TextRange
$tr.Text = "PIWG Action Item XX-XXX"
$firstWord = $tr.Characters(1, 4)
$firstWord.Text = "PIWG"

# --- 2. Remove the stray demo shapes that followed "Problem Background"
$namesToRemove = @(
    "TextBox 2",
    "TextBox 5",
    "TextBox 9",
    "TextBox 10",
    "TextBox 11",
    "TextBox 12",
    "TextBox 13",
    "TextBox 14",
    "TextBox 15",
    "TextBox 16",
    "TextBox 17",
    "TextBox 18",
    "TextBox 19",
    "TextBox 20",
    "TextBox 21",
    "TextBox 22",
    "Rectangle 23",
    "TextBox 24"
)

foreach ($name in $namesToRemove) {
    $s.Shapes.Item($name).Delete()
}
